$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 210.47826
$ws.Range("J9").Value = 846.3333
$ws.Range("L9").Value = 846.3333
$ws.Range("N9").Value = -1184.3333

$ws.Range("H40").Value = 5279.2856
$ws.Range("I40").Value = 7497.5
$ws.Range("J40").Value = 4392
$ws.Range("K40").Value = 7497.5
$ws.Range("L40").Value = 4392
$ws.Range("M40").Value = -7322.5
$ws.Range("N40").Value = -4742

$ws.Range("H62").Value = 18527404
$ws.Range("I62").Value = 27786108
$ws.Range("K62").Value = 27786108
$ws.Range("M62").Value = -27785484

$ws.Range("H63").Value = 29166.666
$ws.Range("I63").Value = 25000
$ws.Range("K63").Value = 25000
$ws.Range("M63").Value = -24376

$ws.Range("H65").Value = 18527404
$ws.Range("I65").Value = 27786108
$ws.Range("K65").Value = 138930540
$ws.Range("M65").Value = -138927420

$ws.Range("H66").Value = 29166.666
$ws.Range("I66").Value = 25000
$ws.Range("K66").Value = 75000
$ws.Range("M66").Value = -71880

$ws.Range("H74").Value = 3710.625
$ws.Range("I74").Value = 3126.5
$ws.Range("K74").Value = 3126.5
$ws.Range("M74").Value = -2190.5

$ws.Range("H77").Value = 3710.625
$ws.Range("I77").Value = 3126.5
$ws.Range("K77").Value = 15632.5
$ws.Range("M77").Value = -10952.5

$ws.Range("H80").Value = 847.9375
$ws.Range("J80").Value = 960.6667
$ws.Range("L80").Value = 2882.0001
$ws.Range("N80").Value = -4878.0001

$ws.Range("H83").Value = 847.9375
$ws.Range("J83").Value = 960.6667
$ws.Range("L83").Value = 8646.0003
$ws.Range("N83").Value = -18630.0003

$ws.Range("H100").Value = 1700.5
$ws.Range("I100").Value = 1901.9
$ws.Range("K100").Value = 1901.9
$ws.Range("M100").Value = -1360.9

$ws.Range("H123").Value = 69622.25
$ws.Range("J123").Value = 69622.25
$ws.Range("L123").Value = 69622.25
$ws.Range("N123").Value = -79422.25

$ws.Range("H138").Value = 3210.1133
$ws.Range("I138").Value = 2099.5
$ws.Range("J138").Value = 3300.7754
$ws.Range("K138").Value = 6298.5
$ws.Range("L138").Value = 9902.3262
$ws.Range("M138").Value = -1158.5
$ws.Range("N138").Value = -20182.3262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6583.892
$ws.Range("I32").Value = 3700.4922
$ws.Range("K32").Value = 3700.4922
$ws.Range("M32").Value = -3413.4922

$ws.Range("H63").Value = 2796.1667
$ws.Range("I63").Value = 2796.1667
$ws.Range("K63").Value = 2796.1667
$ws.Range("M63").Value = -2110.1667

$ws.Range("H66").Value = 2796.1667
$ws.Range("I66").Value = 2796.1667
$ws.Range("K66").Value = 13980.8335
$ws.Range("M66").Value = -10548.8335

$ws.Range("H97").Value = 1958.238
$ws.Range("I97").Value = 1624.9412
$ws.Range("K97").Value = 1624.9412
$ws.Range("M97").Value = -1128.9412

$ws.Range("H110").Value = 4073.923
$ws.Range("I110").Value = 4551
$ws.Range("K110").Value = 4551
$ws.Range("M110").Value = -2506

$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494

$ws.Range("H122").Value = 2726.0667
$ws.Range("I122").Value = 2530.6924
$ws.Range("J122").Value = 3996
$ws.Range("K122").Value = 7592.0772
$ws.Range("L122").Value = 11988
$ws.Range("M122").Value = -5142.0772
$ws.Range("N122").Value = -16888

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 40698
$ws.Range("I82").Value = 17320.8
$ws.Range("J82").Value = 79660
$ws.Range("K82").Value = 17320.8
$ws.Range("L82").Value = 79660
$ws.Range("M82").Value = -16937.8
$ws.Range("N82").Value = -80426

$ws.Range("H85").Value = 40698
$ws.Range("I85").Value = 17320.8
$ws.Range("J85").Value = 79660
$ws.Range("K85").Value = 17320.8
$ws.Range("L85").Value = 79660
$ws.Range("M85").Value = -15994.8
$ws.Range("N85").Value = -82312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3455
$ws.Range("I134").Value = 3069.762
$ws.Range("K134").Value = 9209.286
$ws.Range("M134").Value = -6674.286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 950
$ws.Range("I3").Value = 950
$ws.Range("K3").Value = 2850
$ws.Range("M3").Value = -2738

$ws.Range("H108").Value = 1144.75
$ws.Range("I108").Value = 1144.75
$ws.Range("K108").Value = 3434.25
$ws.Range("M108").Value = -554.25

$ws.Range("H113").Value = 772.6429000000001
$ws.Range("I113").Value = 698.5
$ws.Range("K113").Value = 2095.5
$ws.Range("M113").Value = 74.5

$ws.Range("H131").Value = 1267.2693
$ws.Range("I131").Value = 773.4286
$ws.Range("J131").Value = 1449.2106
$ws.Range("K131").Value = 2320.2858
$ws.Range("L131").Value = 4347.6318
$ws.Range("M131").Value = 2719.7142
$ws.Range("N131").Value = -14427.6318

$ws.Range("H132").Value = 7970.0557
$ws.Range("I132").Value = 15121.125
$ws.Range("J132").Value = 2249.2
$ws.Range("K132").Value = 136090.125
$ws.Range("L132").Value = 20242.8
$ws.Range("M132").Value = -133560.125
$ws.Range("N132").Value = -25302.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3081.3333
$ws.Range("I113").Value = 3006.8333
$ws.Range("J113").Value = 3155.8333
$ws.Range("K113").Value = 3006.8333
$ws.Range("L113").Value = 3155.8333
$ws.Range("M113").Value = -836.8332999999998
$ws.Range("N113").Value = -7495.8333

$ws.Range("H126").Value = 10095.6
$ws.Range("I126").Value = 10095.6
$ws.Range("K126").Value = 30286.8
$ws.Range("M126").Value = -27816.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 7001
$ws.Range("I17").Value = 7001
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 7001
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -6831
$ws.Range("N17").ClearContents()

$ws.Range("H22").Value = 3711
$ws.Range("J22").Value = 5551.6665
$ws.Range("L22").Value = 5551.6665
$ws.Range("N22").Value = -6141.6665

$ws.Range("H27").Value = 3711
$ws.Range("J27").Value = 5551.6665
$ws.Range("L27").Value = 5551.6665
$ws.Range("N27").Value = -5765.6665

$ws.Range("H40").Value = 5793.5
$ws.Range("I40").Value = 4951.4
$ws.Range("K40").Value = 4951.4
$ws.Range("M40").Value = -4815.4

$ws.Range("H61").Value = 4416.2
$ws.Range("I61").Value = 4524
$ws.Range("K61").Value = 4524
$ws.Range("M61").Value = -4322

$ws.Range("H113").Value = 4416.2
$ws.Range("I113").Value = 4524
$ws.Range("K113").Value = 4524
$ws.Range("M113").Value = -2354

$ws.Range("H132").Value = 13497.5
$ws.Range("I132").Value = 12995
$ws.Range("K132").Value = 38985
$ws.Range("M132").Value = -36455

$ws.Range("H136").Value = 6394.4
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H62").Value = 10053.111
$ws.Range("J62").Value = 11663
$ws.Range("L62").Value = 11663
$ws.Range("N62").Value = -12911

$ws.Range("H65").Value = 10053.111
$ws.Range("J65").Value = 11663
$ws.Range("L65").Value = 58315
$ws.Range("N65").Value = -64555

$ws.Range("H113").Value = 1499.2858
$ws.Range("J113").Value = 2002.5
$ws.Range("L113").Value = 6007.5
$ws.Range("N113").Value = -10347.5

$ws.Range("H122").Value = 3382
$ws.Range("I122").Value = 3378.4
$ws.Range("K122").Value = 10135.2
$ws.Range("M122").Value = -7685.200000000001
